$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) "Cropped vs. Uncropped Images" table (around row 29-34): a new "running"
#    placeholder row and a new filled-in result row are inserted, plus a
#    trailing blank spacer row, pushing the next table down by two rows.
#    We splice in three style-s4 rows (copy+insert, which is the only
#    operation that reliably carries the exact - including non-visual -
#    style flags) and then delete one spare row so the net shift is +2,
#    matching the target layout.
# ---------------------------------------------------------------------------
$ws.Rows("31:31").Copy()
$ws.Rows("31:31").Insert()
$excel.CutCopyMode = $false

$ws.Rows("32:32").Copy()
$ws.Rows("33:33").Insert()
$excel.CutCopyMode = $false

$ws.Rows("33:33").Copy()
$ws.Rows("34:34").Insert()
$excel.CutCopyMode = $false

$ws.Rows("35:35").Delete()

# Row 31: new "?" placeholder row. Duplicate the still-original row 15
# (style s=2, the "?" placeholder look) via insert-below + delete-original so
# the exact style (including non-visual flags) is preserved, then set content.
$ws.Rows("15:15").Copy()
$ws.Rows("32:32").Insert()
$excel.CutCopyMode = $false
$ws.Rows("31:31").Delete()

$ws.Cells.Item(31, 1).Value = "Aligned_cropped_CL9_DL1_nobias_200Epoch_0.0001LR1Batch1keep0WD"
$ws.Cells.Item(31, 2).Value = "?"
$ws.Cells.Item(31, 3).Value = "?"

# Row 32 already holds the untouched LR3Batch result (copied down automatically)

# Row 33: new filled-in result row
$ws.Cells.Item(33, 1).Value = "Aligned_CL9_DL1_nobias_200Epoch_0.0001LR1Batch1keep0WD"
$ws.Cells.Item(33, 2).Value = [double]"0.92198325694600003"
$ws.Cells.Item(33, 3).Value = [double]"8.0061481612300006E-2"

# Row 34: blank spacer row (keep style s=4, no content)
$ws.Range("A34:C34").ClearContents()

# ---------------------------------------------------------------------------
# 2) Row 15 (Comparison of batch size / LR1Batch): the "?" placeholder result
#    is now filled in with real numbers. Style changes from s=2 to s=4.
#    Duplicate row 16 (already style s=4) via insert-below + delete-original
#    so row 15 ends up with the exact style, then set its content.
# ---------------------------------------------------------------------------
$ws.Rows("16:16").Copy()
$ws.Rows("17:17").Insert()
$excel.CutCopyMode = $false
$ws.Rows("15:15").Delete()

$ws.Cells.Item(15, 1).Value = "Aligned_CL9_DL1_nobias_200Epoch_0.0001LR1Batch1keep0WD"
$ws.Cells.Item(15, 2).Value = [double]"0.92198325694600003"
$ws.Cells.Item(15, 3).Value = [double]"8.0061481612300006E-2"

# ---------------------------------------------------------------------------
# 3) Sheet view: update the visible top row / selection to match.
# ---------------------------------------------------------------------------
$ws.Application.GoTo($ws.Range("C36"), $true)
$ws.Range("C36").Select()
